# Automated update of EPEX Spot prices workbook
# - "Prix Spot" sheet: add a new day column (AO) with header "24-jul" and its 24 hourly values
# - "Gaz" sheet: append a new row (38) with date 2025-07-22 and its price
# - "CO2" sheet: append a new row (38) with date 2025-07-22 and its price

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Prix Spot": new column AO ("24-jul")
# ---------------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell, formatted like the other day headers (bold / centered / bordered)
$wsPrix.Range("AO1").Value = "24-jul"
$wsPrix.Range("AN1").Copy()
$wsPrix.Range("AO1").PasteSpecial(-4122)

# Hourly values for 24-jul
$prixSpotValues = @{
    2  = 96.23999999999999
    3  = 87.36
    4  = 81.88
    5  = 71.79000000000001
    6  = 63.27
    7  = 80.90000000000001
    8  = 85.62
    9  = 90.44
    10 = 103.04
    11 = 90.83
    12 = 70.90000000000001
    13 = 65.64
    14 = 71.05
    15 = 34.04
    16 = 28
    17 = 37.02
    18 = 26.53
    19 = 41.42
    20 = 74.72
    21 = 92.29000000000001
    22 = 87.93000000000001
    23 = 90.92
    24 = 109.19
    25 = 102.06
}

foreach ($row in $prixSpotValues.Keys) {
    $wsPrix.Cells.Item($row, 41).Value = $prixSpotValues[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Gaz": append row 38 (2025-07-22 / 32.625)
# ---------------------------------------------------------------------------
# The date column is stored as plain text in this workbook (not a real Excel
# date), so force a text number format before writing the value to stop
# Excel from auto-converting "2025-07-22" into a date serial number, then
# restore the plain "Normal" style so the cell keeps the same (style-less)
# look as the rest of the column.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A38").NumberFormat = "@"
$wsGaz.Range("A38").Value = "2025-07-22"
$wsGaz.Range("A38").Style = "Normal"
$wsGaz.Range("B38").Value = 32.625

# ---------------------------------------------------------------------------
# Sheet "CO2": append row 38 (2025-07-22 / 68.25)
# ---------------------------------------------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A38").NumberFormat = "@"
$wsCo2.Range("A38").Value = "2025-07-22"
$wsCo2.Range("A38").Style = "Normal"
$wsCo2.Range("B38").Value = 68.25

Write-Output "done"
